$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 368.57144
$ws.Range("I2").Value = 400
$ws.Range("K2").Value = 400
$ws.Range("M2").Value = -287
$ws.Range("H19").Value = 849428.9399999999
$ws.Range("I19").Value = 1385479.8
$ws.Range("K19").Value = 1385479.8
$ws.Range("M19").Value = -1385304.8
$ws.Range("H28").Value = 855.2222
$ws.Range("I28").Value = 855.2222
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 855.2222
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -370.2222
$ws.Range("N28").ClearContents()
$ws.Range("H98").Value = 8619.048000000001
$ws.Range("J98").Value = 9818
$ws.Range("L98").Value = 9818
$ws.Range("N98").Value = -12814
$ws.Range("H107").Value = 1293.7894
$ws.Range("I107").Value = 1487.7778
$ws.Range("J107").Value = 1119.2
$ws.Range("K107").Value = 1487.7778
$ws.Range("L107").Value = 1119.2
$ws.Range("M107").Value = 432.2221999999999
$ws.Range("N107").Value = -4959.2
$ws.Range("H116").Value = 843784
$ws.Range("I116").Value = 5001895
$ws.Range("J116").Value = 12161.8
$ws.Range("K116").Value = 5001895
$ws.Range("L116").Value = 12161.8
$ws.Range("M116").Value = -4998453
$ws.Range("N116").Value = -19045.8
$ws.Range("H122").Value = 8619.048000000001
$ws.Range("J122").Value = 9818
$ws.Range("L122").Value = 29454
$ws.Range("N122").Value = -34354
$ws.Range("H132").Value = 29416592
$ws.Range("I132").Value = 33338166
$ws.Range("J132").Value = 4776.75
$ws.Range("K132").Value = 100014498
$ws.Range("L132").Value = 14330.25
$ws.Range("M132").Value = -100011968
$ws.Range("N132").Value = -19390.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 974.1429000000001
$ws.Range("I2").Value = 863.8
$ws.Range("J2").Value = 1250
$ws.Range("K2").Value = 863.8
$ws.Range("L2").Value = 1250
$ws.Range("M2").Value = -750.8
$ws.Range("N2").Value = -1476
$ws.Range("H63").Value = 11545941
$ws.Range("I63").Value = 15392588
$ws.Range("K63").Value = 15392588
$ws.Range("M63").Value = -15391902
$ws.Range("H66").Value = 11545941
$ws.Range("I66").Value = 15392588
$ws.Range("K66").Value = 76962940
$ws.Range("M66").Value = -76959508
$ws.Range("H116").Value = 974.1429000000001
$ws.Range("I116").Value = 863.8
$ws.Range("J116").Value = 1250
$ws.Range("K116").Value = 863.8
$ws.Range("L116").Value = 1250
$ws.Range("M116").Value = 1430.2
$ws.Range("N116").Value = -5838
$ws.Range("H139").Value = 41593
$ws.Range("J139").Value = 41593
$ws.Range("L139").Value = 41593
$ws.Range("N139").Value = -51873

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 974.1429000000001
$ws.Range("I3").Value = 863.8
$ws.Range("J3").Value = 1250
$ws.Range("K3").Value = 863.8
$ws.Range("L3").Value = 1250
$ws.Range("M3").Value = -749.8
$ws.Range("N3").Value = -1478
$ws.Range("H97").Value = 35277.43
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 35277.43
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 35277.43
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -37259.43
$ws.Range("H138").Value = 40883.332
$ws.Range("J138").Value = 40883.332
$ws.Range("L138").Value = 40883.332
$ws.Range("N138").Value = -51163.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7788.2104
$ws.Range("I31").Value = 1494.5
$ws.Range("K31").Value = 1494.5
$ws.Range("M31").Value = -1199.5
$ws.Range("H34").Value = 7788.2104
$ws.Range("I34").Value = 1494.5
$ws.Range("K34").Value = 1494.5
$ws.Range("M34").Value = -1292.5
$ws.Range("H134").Value = 4035.3416
$ws.Range("I134").Value = 4351.2334
$ws.Range("J134").Value = 3173.818
$ws.Range("K134").Value = 13053.7002
$ws.Range("L134").Value = 9521.454000000002
$ws.Range("M134").Value = -10518.7002
$ws.Range("N134").Value = -14591.454
$ws.Range("H138").Value = 42903.332
$ws.Range("J138").Value = 42903.332
$ws.Range("L138").Value = 42903.332
$ws.Range("N138").Value = -53183.332
$ws.Range("H140").Value = 92187.27
$ws.Range("J140").Value = 92187.27
$ws.Range("L140").Value = 92187.27
$ws.Range("N140").Value = -102547.27
$ws.Range("H141").Value = 30842.6
$ws.Range("J141").Value = 30842.6
$ws.Range("L141").Value = 30842.6
$ws.Range("N141").Value = -41202.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 2778457.8
$ws.Range("I113").Value = 573.78125
$ws.Range("J113").Value = 9616326
$ws.Range("K113").Value = 1721.34375
$ws.Range("L113").Value = 28848978
$ws.Range("M113").Value = 448.65625
$ws.Range("N113").Value = -28853318
$ws.Range("H131").Value = 815.6161499999999
$ws.Range("I131").Value = 295.25
$ws.Range("J131").Value = 837.5263
$ws.Range("K131").Value = 885.75
$ws.Range("L131").Value = 2512.5789
$ws.Range("M131").Value = 4154.25
$ws.Range("N131").Value = -12592.5789

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H140").Value = 38492.777
$ws.Range("J140").Value = 38492.777
$ws.Range("L140").Value = 38492.777
$ws.Range("N140").Value = -48852.777

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H87").Value = 50000
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()
$ws.Range("H90").Value = 50000
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()
$ws.Range("H139").Value = 45840
$ws.Range("J139").Value = 45840
$ws.Range("L139").Value = 45840
$ws.Range("N139").Value = -56120
$ws.Range("H140").Value = 67882.836
$ws.Range("J140").Value = 67882.836
$ws.Range("L140").Value = 67882.836
$ws.Range("N140").Value = -78242.836
$ws.Range("H141").Value = 31894.564
$ws.Range("J141").Value = 31894.564
$ws.Range("L141").Value = 31894.564
$ws.Range("N141").Value = -42254.564

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 62903.453
$ws.Range("J46").Value = 62903.453
$ws.Range("L46").Value = 62903.453
$ws.Range("N46").Value = -63365.453
$ws.Range("H122").Value = 4888.1665
$ws.Range("I122").Value = 3677.6428
$ws.Range("J122").Value = 9125
$ws.Range("K122").Value = 11032.9284
$ws.Range("L122").Value = 27375
$ws.Range("M122").Value = -8582.928400000001
$ws.Range("N122").Value = -32275
$ws.Range("H132").Value = 30306984
$ws.Range("I132").Value = 2140.6667
$ws.Range("K132").Value = 6422.000100000001
$ws.Range("M132").Value = -3892.000100000001
$ws.Range("H134").Value = 62903.453
$ws.Range("J134").Value = 62903.453
$ws.Range("L134").Value = 188710.359
$ws.Range("N134").Value = -193780.359
$ws.Range("H138").Value = 42999.5
$ws.Range("J138").Value = 42999.5
$ws.Range("L138").Value = 42999.5
$ws.Range("N138").Value = -53279.5
$ws.Range("H139").Value = 41260.555
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 41260.555
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 41260.555
$ws.Range("M139").ClearContents()
$ws.Range("N139").Value = -51540.555
$ws.Range("H140").Value = 30761.2
$ws.Range("J140").Value = 30761.2
$ws.Range("L140").Value = 30761.2
$ws.Range("N140").Value = -41121.2
$ws.Range("H141").Value = 31315
$ws.Range("J141").Value = 31315
$ws.Range("L141").Value = 31315
$ws.Range("N141").Value = -41675
